# Auto commit at 2025-11-02  8:56:44.70
# Refreshes the daily "Metrics" feed values and the derived numbers on the
# "today" sheet that are pasted (not formula-linked) for the first four
# metric rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metrics sheet: raw feed values for B2:B13 refreshed for the new day.
# ---------------------------------------------------------------------
$wsMetrics = $wb.Worksheets.Item("Metrics")

$wsMetrics.Range("B2").Value  = 12362.01
$wsMetrics.Range("B3").Value  = 10742.68
$wsMetrics.Range("B4").Value  = 3777.95
$wsMetrics.Range("B5").Value  = 540
$wsMetrics.Range("B6").Value  = 4808607.76
$wsMetrics.Range("B7").Value  = 4052819.3600000003
$wsMetrics.Range("B8").Value  = 1410737.78
$wsMetrics.Range("B9").Value  = 186747
$wsMetrics.Range("B10").Value = 33273988.750000004
$wsMetrics.Range("B11").Value = 31369037.210000001
$wsMetrics.Range("B12").Value = 11692459.82
$wsMetrics.Range("B13").Value = 1284377

# ---------------------------------------------------------------------
# 2. "today" sheet: rows 11-14 (B/E/F) get the same refreshed numbers
#    pasted as static values (this overwrites/removes their old
#    formulas, matching how the automation snapshot was produced).
#    Rows 15-22 keep their live formulas - they recompute on their own
#    once the Metrics sheet above changes.
# ---------------------------------------------------------------------
$wsToday = $wb.Worksheets.Item("today")

$wsToday.Range("B11").Value = 12362.01
$wsToday.Range("E11").Value = 12362.01
$wsToday.Range("F11").Value = 12362.01

$wsToday.Range("B12").Value = 10742.68
$wsToday.Range("E12").Value = 10742.68
$wsToday.Range("F12").Value = 10742.68

$wsToday.Range("B13").Value = 3777.95
$wsToday.Range("E13").Value = 3777.95
$wsToday.Range("F13").Value = 3777.95

$wsToday.Range("B14").Value = 540
$wsToday.Range("E14").Value = 540
$wsToday.Range("F14").Value = 540

# ---------------------------------------------------------------------
# 3. New (empty, but styled) helper cells in column J for rows 19-22,
#    matching the style already used by the neighbouring E/F/I cells.
#    This also brings columns H and J into the sheet's <cols> list.
# ---------------------------------------------------------------------
$wsToday.Range("J19").NumberFormat = "#,##0.00_ "
$wsToday.Range("J20").NumberFormat = "#,##0.00_ "
$wsToday.Range("J21").NumberFormat = "#,##0.00_ "
$wsToday.Range("J22").NumberFormat = "#,##0.00_ "

$wsToday.Columns.Item(8).ColumnWidth  = 17.643
$wsToday.Columns.Item(10).ColumnWidth = 19.501

# ---------------------------------------------------------------------
# 4. Restore the cursor/selection on each sheet the way it was left in
#    the saved workbook (Metrics selection first, "today" last so it
#    stays the active sheet/tab).
# ---------------------------------------------------------------------
$wsMetrics.Range("E13").Select()
$wsToday.Range("I14").Select()

"done"
